# Data retrieved - Wed Jul 14 18:34:50 UTC 2021
#
# Apply a tiny correction to the timestamp of the last existing row (77),
# then append a new data row (78) with the latest job-number snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 77: timestamp value refined (floating point re-save of the date/time).
$ws.Range("A77").Value = 44390.76765445602

# Row 78: newly retrieved data point.
$ws.Range("A78").Value = 44391.77418872772
$ws.Range("B78").Value = 80254
$ws.Range("C78").Value = 67613
$ws.Range("D78").Value = 3598
$ws.Range("E78").Value = 2247
$ws.Range("F78").Value = 1612
$ws.Range("G78").Value = 21327
$ws.Range("H78").Value = 1544
$ws.Range("I78").Value = 914
$ws.Range("J78").Value = 197
